$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header labels for the new columns (AD, AE, AF), matching style of existing headers (A1:AC1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style from an existing header cell (AC1) to the new header cells
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Fill in season record values for each data row (2 through 43)
for ($row = 2; $row -le 43; $row++) {
    $ws.Cells.Item($row, 30).Value = 84  # AD = Wins
    $ws.Cells.Item($row, 31).Value = 78  # AE = Losses
    $ws.Cells.Item($row, 32).Value = 0   # AF = Ties
}
